# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.741.79'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.11%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.033.46'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '228.01'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.95%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.608'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '60.07'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.49%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -1.25%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '2.333.63'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.28%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.04'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.25%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.769'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("E16").Value = '  -2.88%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.019.20'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.77%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.723.27'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.05%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '69.59'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  -5.97%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0824'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.97%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '223.62'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("E23").Value = '  +0.14%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.40'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("E25").Value = '  +3.16%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '167.34'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +8.26%  '
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("E34").Value = '  +0.24%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.50'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.52%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.45'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.51%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.05%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.43'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.38%  '
$ws.Range("E39").Value = '  +0.01%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.12'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +7.10%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.536.28'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.29%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0216'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.75%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '96.28'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.80'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.58%  '
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  -0.98%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.23%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.10'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.58%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.223.27'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.35%  '
